# Stock_Records.xlsx update
#
# Simulates a stock-price-tracker run: today's placeholder sheet ("Sheet6")
# becomes tomorrow's (now-empty) placeholder sheet ("Sheet"), the previous
# run's data sheet ("Sheet") is archived as "Sheet7", and a brand-new sheet
# ("Sheet8") is appended holding this run's freshly fetched ticker/price
# data -- including the known bug where the value used to break out of the
# input loop gets appended to the ticker list too.

function Set-TextValue($range, [string]$value) {
    # Force literal text storage (matches the source file's t="inlineStr"
    # cells) even for values that look numeric/date-like, the way typing an
    # apostrophe-prefixed entry into Excel does. ClearFormats() afterwards
    # drops the transient "quote prefix" cell style so no stray style is
    # left behind on the cell.
    $range.Value = "'" + $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Sheet 2 ("Sheet") must be renamed out of the way first, since Sheet 1
#     ("Sheet6") is about to claim the name "Sheet". ---
$wsArchive = $wb.Worksheets.Item(2)
$wsArchive.Name = "Sheet7"

$wsToday = $wb.Worksheets.Item(1)
$wsToday.Name = "Sheet"

# --- "Sheet" (was "Sheet6"): drop yesterday's ticker list, bump the date,
#     leave just the header row behind. ---
$wsToday.Range("C2:D6").ClearContents()
Set-TextValue $wsToday.Range("A1") "01/05/21"
$wsToday.Range("A1").Select() | Out-Null

# --- "Sheet7" (was "Sheet"): this is last run's completed data; bump its
#     date forward. ---
Set-TextValue $wsArchive.Range("A1") "01/12/21"
$wsArchive.Range("A1").Select() | Out-Null

# --- New "Sheet8": this run's freshly collected prices. Reproduces the
#     known bug -- "1" (whatever was typed to break out of the input loop)
#     ends up appended as a 4th "ticker" with no matching price. ---
$wsNew = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "Sheet8"

$wsNew.PageSetup.LeftMargin = 54
$wsNew.PageSetup.RightMargin = 54
$wsNew.PageSetup.TopMargin = 72
$wsNew.PageSetup.BottomMargin = 72
$wsNew.PageSetup.HeaderMargin = 36
$wsNew.PageSetup.FooterMargin = 36

Set-TextValue $wsNew.Range("A1") "01/12/21"
Set-TextValue $wsNew.Range("C1") "Ticker"
Set-TextValue $wsNew.Range("D1") "Price"

Set-TextValue $wsNew.Range("C2") "tsla"
Set-TextValue $wsNew.Range("D2") "849.44"

Set-TextValue $wsNew.Range("C3") "ge"
Set-TextValue $wsNew.Range("D3") "11.78"

Set-TextValue $wsNew.Range("C4") "gme"
Set-TextValue $wsNew.Range("D4") "19.95"

Set-TextValue $wsNew.Range("C5") "1"

$wsNew.Range("A1").Select() | Out-Null

# "Sheet7" is the active/selected tab (not the brand-new "Sheet8").
$wsArchive.Activate() | Out-Null
